# Consolidate player/team statistics after unifying teams with name changes.
# Columns V (RECUPEROS) and W (PERDIDAS), plus the LOCAL/VISITANTE splits
# BQ/BR (RECUPEROS) and BS/BT (PERDIDAS), are recomputed from per-game
# averages into season totals (value * PJ games played) for rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ V = 114; W = 201; BQ = 57;  BR = 57;  BS = 111; BT = 90  }
    3  = @{ V = 86;  W = 181; BQ = 39;  BR = 47;  BS = 94;  BT = 87  }
    4  = @{ V = 89;  W = 198; BQ = 45;  BR = 44;  BS = 135; BT = 63  }
    5  = @{ V = 98;  W = 175; BQ = 51;  BR = 47;  BS = 90;  BT = 85  }
    6  = @{ V = 107; W = 162; BQ = 57;  BR = 50;  BS = 74;  BT = 88  }
    7  = @{ V = 121; W = 139; BQ = 74;  BR = 47;  BS = 76;  BT = 63  }
    8  = @{ V = 117; W = 150; BQ = 26;  BR = 91;  BS = 34;  BT = 116 }
    9  = @{ V = 93;  W = 160; BQ = 47;  BR = 46;  BS = 82;  BT = 78  }
    10 = @{ V = 141; W = 207; BQ = 80;  BR = 61;  BS = 93;  BT = 114 }
    11 = @{ V = 97;  W = 190; BQ = 49;  BR = 48;  BS = 94;  BT = 96  }
    12 = @{ V = 95;  W = 171; BQ = 46;  BR = 49;  BS = 92;  BT = 79  }
    13 = @{ V = 73;  W = 218; BQ = 31;  BR = 42;  BS = 110; BT = 108 }
    14 = @{ V = 115; W = 173; BQ = 68;  BR = 47;  BS = 107; BT = 66  }
    15 = @{ V = 102; W = 206; BQ = 50;  BR = 52;  BS = 111; BT = 95  }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
